$wb = $excel.ActiveWorkbook

# Sheet "y_fitted_on_begin_2016": column B (y_value), rows 2-38 (years 1980-2016)
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$sheet1Values = @(
    22.64913599676362,
    26.07962720512022,
    25.40006367834748,
    22.33393874347184,
    21.60078926660292,
    21.54596395854735,
    20.40716118038087,
    22.6178837310137,
    23.92411380949441,
    25.68349066168738,
    24.99579767562324,
    24.20210259809138,
    25.13618058040023,
    24.15384176160572,
    24.26194429458621,
    23.5123693637873,
    23.83121168124817,
    23.10755366924378,
    22.20881297484215,
    22.68303995050423,
    22.82747497702729,
    25.49153903092111,
    23.30883597945532,
    24.49711985839617,
    24.52169385200925,
    26.60455794703042,
    26.56685880142078,
    26.12554169201768,
    23.68750514529793,
    20.82330250365575,
    20.54697834972677,
    20.98980037393388,
    19.93571157408378,
    19.52361731089949,
    18.21546479360531,
    18.5352084992968,
    18.61445526537136
)
for ($i = 0; $i -lt $sheet1Values.Count; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $sheet1Values[$i]
}

# Sheet "y_fitted_on_begin_2021": column B (y_value), rows 2-43 (years 1980-2021)
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$sheet3Values = @(
    22.84260586636905,
    26.01623118361493,
    25.35453124045803,
    22.36900506786568,
    21.65512778037565,
    21.60174365689832,
    20.49287641506989,
    22.64548603647859,
    23.91737944250828,
    25.63050782156707,
    24.96089211698524,
    24.18806078187685,
    25.09758479758315,
    24.14106857069841,
    24.24632942874948,
    23.51645846075188,
    23.8269194184519,
    23.1222840951255,
    22.24716846138544,
    22.70892950428498,
    22.84956778911767,
    25.44360199427803,
    23.31827532770089,
    24.47532296809371,
    24.49925098885453,
    26.52736315053029,
    26.49065499665772,
    26.06093872334181,
    23.68699047447056,
    20.89807869433864,
    20.6290182316384,
    21.06019986029055,
    20.03381973324703,
    19.63255812961255,
    18.35879281962764,
    18.67013147090224,
    18.7472950892264,
    18.24665749203839,
    19.17292881674655,
    21.02786746990155,
    19.90483625840966,
    20.35469851372361
)
for ($i = 0; $i -lt $sheet3Values.Count; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 2).Value = $sheet3Values[$i]
}
